$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.710.79'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -3.00%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.097.18'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.97%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.16%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '345.53'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.98%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.009'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.15%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5147'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.36%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4399'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.42%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09275'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.64%  '

# Row 10
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.52'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.56%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.173'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.85%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.87'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.12%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.087.51'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.97%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.285'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.97%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.762'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.40%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '99.58'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.52%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001154'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.54%  '

# Row 18
$ws.Range('E18').Value = '  -0.24%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '20.88'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +7.06%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06646'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.12%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.008'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.22%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.205'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.18%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.752.43'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.11%  '

# Row 24
$ws.Range('E24').Value = '  -1.44%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.321'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.01%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.345.98'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.62%  '

# Row 27
$ws.Range('E27').Value = '  -2.72%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.528'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.50%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '161.93'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.76%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.12'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.86%  '

# Row 31
$ws.Range('E31').Value = '  -6.99%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1051'
$ws.Range('D32').ClearFormats()

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.655'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.81%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.173'
$ws.Range('D34').ClearFormats()

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.939'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.90%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.172'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.50%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.32'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.47%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02577'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.81%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06722'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.83%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.47'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.22%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6861'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.64%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2222'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.61%  '

# Row 43
$ws.Range('E43').Value = '  +1.82%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6634'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.73%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.27'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.75%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.321'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.19%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.629'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.42%  '

# Row 48
$ws.Range('E48').Value = '  -5.88%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.220'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.84%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '82.23'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.09%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3316'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.71%  '
